# Atualização de bases das ligas, do dia: 19-04-2024 às 00:38
# Adds two new match rows (134 and 135) to the "India Super League" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (133) into the
# two new rows, so the new cells pick up the same cell styles (bold/border
# id column, date-time formatted column) without creating extra style
# definitions. Only the columns that will actually contain data are copied
# (A:G and K:AA) -- columns H,I,J,AB,AC stay empty on these new rows, same
# as in the source data.
$ws.Range("A133:G133").Copy()
$ws.Range("A134:G135").PasteSpecial(-4122)
$ws.Range("K133:AA133").Copy()
$ws.Range("K134:AA135").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 134 - Odisha FC vs Kerala Blasters
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 8103573
$ws.Range("C134").Value = "India Super League"
$ws.Range("D134").Value = "India Super League"
$ws.Range("E134").Value = 45401.45833333334
$ws.Range("F134").Value = "Odisha FC"
$ws.Range("G134").Value = "Kerala Blasters"
$ws.Range("K134").Value = 1.65
$ws.Range("L134").Value = 3.7
$ws.Range("M134").Value = 4.5
$ws.Range("N134").Value = 1.65
$ws.Range("O134").Value = 3.75
$ws.Range("P134").Value = 4.5
$ws.Range("Q134").Value = -0.75
$ws.Range("R134").Value = 1.85
$ws.Range("S134").Value = 1.95
$ws.Range("T134").Value = 2.75
$ws.Range("U134").Value = 1.8
$ws.Range("V134").Value = 2
$ws.Range("W134").Value = 0
$ws.Range("X134").Value = 0
$ws.Range("Y134").Value = 0
$ws.Range("Z134").Value = 0
$ws.Range("AA134").Value = 0

# Row 135 - FC Goa vs Chennaiyin FC
$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 8103574
$ws.Range("C135").Value = "India Super League"
$ws.Range("D135").Value = "India Super League"
$ws.Range("E135").Value = 45402.45833333334
$ws.Range("F135").Value = "FC Goa"
$ws.Range("G135").Value = "Chennaiyin FC"
$ws.Range("K135").Value = 1.4
$ws.Range("L135").Value = 4.5
$ws.Range("M135").Value = 6
$ws.Range("N135").Value = 1.4
$ws.Range("O135").Value = 4.5
$ws.Range("P135").Value = 5.75
$ws.Range("Q135").Value = -1.25
$ws.Range("R135").Value = 1.9
$ws.Range("S135").Value = 1.9
$ws.Range("T135").Value = 3
$ws.Range("U135").Value = 1.8
$ws.Range("V135").Value = 2
$ws.Range("W135").Value = 0
$ws.Range("X135").Value = 0
$ws.Range("Y135").Value = 0
$ws.Range("Z135").Value = 0
$ws.Range("AA135").Value = 0
